# Updated cryptos list - refresh Price / Volume(1h) figures, and fix the
# swapped dogwifhat/Cosmos rows (49 & 50).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column holds plain-text numbers (e.g. "1.00", "0.0000120") in
# the source file. Force the column to Text format first so Excel's
# automatic type inference does not strip meaningful trailing zeros when
# the new values are written.
$ws.Range("D2:D51").NumberFormat = "@"

function Set-Row {
    param([int]$Row, [string]$Price, [string]$Volume)
    if ($Price -ne "") {
        $ws.Range("D$Row").Value = $Price
    }
    if ($Volume -ne "") {
        $ws.Range("E$Row").Value = "  $Volume  "
    }
}

Set-Row 2  "66.835.14"  "-3.17%"
Set-Row 3  "3.476.94"   "-2.90%"
Set-Row 4  ""           "+0.26%"
Set-Row 5  "604.63"     "-3.17%"
Set-Row 6  "148.61"     "-5.56%"
Set-Row 7  "3.473.82"   "-2.92%"
Set-Row 8  ""           "+0.00%"
Set-Row 9  ""           "-2.05%"
Set-Row 10 ""           "-3.54%"
Set-Row 11 "7.54"       "+2.24%"
Set-Row 12 ""           "-3.68%"
Set-Row 13 ""           "-4.25%"
Set-Row 14 "31.81"      "-5.15%"
Set-Row 15 "4.064.06"   "-2.60%"
Set-Row 16 "3.476.05"   "-2.59%"
Set-Row 17 "66.899.60"  "-3.46%"
Set-Row 18 ""           "-0.55%"
Set-Row 19 ""           "-5.74%"
Set-Row 20 ""           "-4.23%"
Set-Row 21 ""           "-0.23%"
Set-Row 22 "440.12"     "-4.73%"
Set-Row 23 ""           "-5.11%"
Set-Row 24 "79.33"      ""
Set-Row 25 ""           "+0.12%"
Set-Row 26 "3.612.21"   "-2.76%"
Set-Row 27 "0.0000120"  "-9.19%"
Set-Row 28 ""           "-7.93%"
Set-Row 29 "8.41"       "-7.70%"
Set-Row 30 ""           "-3.84%"
Set-Row 31 ""           "-6.40%"
Set-Row 32 ""           "-1.62%"
Set-Row 33 ""           "-0.18%"
Set-Row 34 "25.47"      "-3.54%"
Set-Row 35 ""           "-6.44%"
Set-Row 36 "3.465.80"   "-2.82%"
Set-Row 37 ""           "-6.77%"
Set-Row 38 "7.95"       "-5.01%"
Set-Row 39 ""           "-0.01%"
Set-Row 40 "1.00"       "+0.33%"
Set-Row 41 "176.70"     "-1.77%"
Set-Row 42 "0.0892"     "-3.34%"
Set-Row 43 "2.14"       "-11.12%"
Set-Row 44 "5.44"       "-3.83%"
Set-Row 45 ""           "-1.77%"
Set-Row 46 "29.38"      "-4.82%"
Set-Row 47 "46.32"      "+0.93%"
Set-Row 48 ""           "-9.14%"

# Rows 49 & 50 swap places: Cosmos moves up to rank 47 (row 49),
# dogwifhat moves down to rank 48 (row 50); both get refreshed figures.
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "7.49"
$ws.Range("E49").Value = "  -4.53%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").Value = "2.46"
$ws.Range("E50").Value = "  -8.66%  "

Set-Row 51 "0.987"      "-4.61%"
